$d = $word.ActiveDocument

# The existing "_GoBack" bookmark (left over from the previous edit session,
# sitting on the "String raza;" paragraph) no longer reflects where the user
# last typed, so Word drops it from there...
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ...and re-creates it at the new editing location: the user added a blank
# line after "int numeroVidas;" and then typed a new line of text, so that's
# where the cursor (and therefore the new "_GoBack" mark) ends up.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>añlkjdlkjfgñlkalkjlkñfjaslkjfdlkj</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')
